$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1: internal/key names ---
$ws.Range("J1").Value = "maxhp"
$ws.Range("K1").Value = "maxhungry"
$ws.Range("L1").Value = "strength"
$ws.Range("M1").Value = "magic"
$ws.Range("N1").Value = "speed"
$ws.Range("O1").Value = "mobility"
$ws.Range("P1").Value = "energy"
$ws.Range("Q1").Value = "taunt"
$ws.Range("R1").Value = "go_ahead"

# --- Row 2: Simplified Chinese display names ---
$ws.Range("J2").Value = "血量"
$ws.Range("K2").Value = "饥饿度"
$ws.Range("L2").Value = "力量"
$ws.Range("M2").Value = "法力"
$ws.Range("N2").Value = "速度"
$ws.Range("O2").Value = "行动力"
$ws.Range("P2").Value = "精力"
$ws.Range("Q2").Value = "嘲讽值"
$ws.Range("R2").Value = "出发"

# --- Row 3: English display names ---
$ws.Range("J3").Value = "MaxHP"
$ws.Range("K3").Value = "MaxHungry"
$ws.Range("L3").Value = "STR"
$ws.Range("M3").Value = "MAG"
$ws.Range("N3").Value = "SPD"
$ws.Range("O3").Value = "DEX"
$ws.Range("P3").Value = "CON"
$ws.Range("Q3").Value = "TAU"
$ws.Range("R3").Value = "Let's Go"

# --- Window/selection state: scroll so column F is leftmost, active cell Q4 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Q4").Select()
